# Backup before the s_RA reformulation: update Sheet1!A1 with the new value.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")
$ws.Range("A1").Value = 5.5061999999999971
